# "add invoice column in download file"
#
# Inserts a new "Invoice Id" column right after "Customer Name" (i.e. the
# new column becomes column B), pushing every existing column one slot to
# the right. The header row gets the literal label "Invoice Id" and the
# data row gets the merge placeholder "{spare:purchase_invoice_id}",
# matching the style already used by their respective rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from column B onward to the right, freeing up column B.
$ws.Columns.Item(2).Insert()

# New header cell (row 1) + new data/template cell (row 2).
$ws.Range("B1").Value = "Invoice Id"
$ws.Range("B2").Value = "{spare:purchase_invoice_id}"

# Match formatting of the other cells in each row.
$ws.Range("A2:M2").Font.Name = "Cambria"
$ws.Range("A2:M2").HorizontalAlignment = -4131
$ws.Range("A2:M2").VerticalAlignment = -4107

$ws.Range("A1:M1").HorizontalAlignment = -4108
$ws.Range("A1:M1").VerticalAlignment = -4107
